$wb = $excel.ActiveWorkbook

# Sheet "展览": row 4 -> F4 964 -> 968, row 6 -> F6 56 -> 58
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F4").Value = 968
$wsExhibit.Range("F6").Value = 58

# Sheet "全部类型": row 5 -> F5 964 -> 968, row 7 -> F7 56 -> 58
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F5").Value = 968
$wsAll.Range("F7").Value = 58
